# Generate Report for Handoff
# Adds a new row for file "a7b876d7-12d5-4694-aa5f-3a4ae396451eo...md" that is
# now "Ready for handoff" across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$newFileName = "a7b876d7-12d5-4694-aa5f-3a4ae396451eooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newPathAndName = "e2e\" + $newFileName
$newHyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a31b96908940a603aa194cca95ce81895a57151/e2e/" + $newFileName
$newZhCnXlf = "a7b876d7-12d5-4694-aa5f-3a4ae396451eoooooooooooooooooooooooooooooooooooooooo.058a10becfaa90a64ebbb2623eb1e69327bd0d30.zh-cn.xlf"
$newDeDeXlf = "a7b876d7-12d5-4694-aa5f-3a4ae396451eoooooooooooooooooooooooooooooooooooooooo.058a10becfaa90a64ebbb2623eb1e69327bd0d30.de-de.xlf"
$readyStatus = "Ready for handoff"
$handoffDateTime = "2016-08-28 12:27:55"
$zhCnHandoffDateTime = "2016-08-28 12:27:51"
$deDeHandoffDateTime = "2016-08-28 12:27:55"
$zeroDateTime = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet: append row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newFileName
$wsOverview.Range("B3").Value = $newPathAndName
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $readyStatus
$wsOverview.Range("F3").Value = $readyStatus
$wsOverview.Range("G3").Value = $handoffDateTime
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newHyperlinkUrl, "", "", $newPathAndName) | Out-Null
$wsOverview.Range("B3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# zh-cn sheet: append row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $newFileName
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $readyStatus
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $newZhCnXlf
$wsZhCn.Range("H3").Value = $zhCnHandoffDateTime
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = $zeroDateTime
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newHyperlinkUrl, "", "", $newFileName) | Out-Null
$wsZhCn.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# de-de sheet: append row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $newFileName
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $readyStatus
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $newDeDeXlf
$wsDeDe.Range("H3").Value = $deDeHandoffDateTime
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = $zeroDateTime
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newHyperlinkUrl, "", "", $newFileName) | Out-Null
$wsDeDe.Range("A3").Style = "HyperLink"
